# Scheduled-runner market data refresh for the Adamantoise_Profits workbook.
# Updates currentAveragePrice(NQ/HQ)/LevePrice(NQ/HQ)/LeveProfit(NQ/HQ) columns (H:N)
# for each Disciple of the Hand sheet, mirroring the upstream API pull.
$wb = $excel.ActiveWorkbook

# ----- ALC (38 cell updates) -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 27778738
$ws.Range("I11").Value = 27778738
$ws.Range("K11").Value = 27778738
$ws.Range("M11").Value = -27778598
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H106").Value = 8337118
$ws.Range("I106").Value = 8337118
$ws.Range("K106").Value = 8337118
$ws.Range("M106").Value = -8336487
$ws.Range("H113").Value = 1749.5
$ws.Range("I113").Value = 1672
$ws.Range("J113").Value = 1920
$ws.Range("K113").Value = 1672
$ws.Range("L113").Value = 1920
$ws.Range("M113").Value = 1582
$ws.Range("N113").Value = -8428
$ws.Range("H116").Value = 35529.6
$ws.Range("I116").Value = 40662
$ws.Range("K116").Value = 40662
$ws.Range("M116").Value = -37220
$ws.Range("H117").Value = 109087.8
$ws.Range("J117").Value = 109087.8
$ws.Range("L117").Value = 109087.8
$ws.Range("N117").Value = -118265.8
$ws.Range("H138").Value = 2675.1746
$ws.Range("I138").Value = 1851.4231
$ws.Range("J138").Value = 3254.027
$ws.Range("K138").Value = 5554.2693
$ws.Range("L138").Value = 9762.081
$ws.Range("M138").Value = -414.2692999999999
$ws.Range("N138").Value = -20042.081

# ----- ARM (34 cell updates) -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3039.2354
$ws.Range("I61").Value = 3243.2222
$ws.Range("K61").Value = 3243.2222
$ws.Range("M61").Value = -3031.2222
$ws.Range("H63").Value = 3641.6667
$ws.Range("I63").Value = 2314.6667
$ws.Range("J63").Value = 4968.6665
$ws.Range("K63").Value = 2314.6667
$ws.Range("L63").Value = 4968.6665
$ws.Range("M63").Value = -1628.6667
$ws.Range("N63").Value = -6340.6665
$ws.Range("H66").Value = 3641.6667
$ws.Range("I66").Value = 2314.6667
$ws.Range("J66").Value = 4968.6665
$ws.Range("K66").Value = 11573.3335
$ws.Range("L66").Value = 24843.3325
$ws.Range("M66").Value = -8141.333500000001
$ws.Range("N66").Value = -31707.3325
$ws.Range("H74").Value = 2411.25
$ws.Range("I74").Value = 2172
$ws.Range("K74").Value = 2172
$ws.Range("M74").Value = -1298
$ws.Range("H77").Value = 2411.25
$ws.Range("I77").Value = 2172
$ws.Range("K77").Value = 10860
$ws.Range("M77").Value = -6492
$ws.Range("H132").Value = 2679.8445
$ws.Range("I132").Value = 2377.2354
$ws.Range("K132").Value = 7131.706200000001
$ws.Range("M132").Value = -4601.706200000001
$ws.Range("H136").Value = 3039.2354
$ws.Range("I136").Value = 3243.2222
$ws.Range("K136").Value = 9729.6666
$ws.Range("M136").Value = -7179.6666

# ----- BSM (25 cell updates) -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1525.6
$ws.Range("I86").Value = 5028
$ws.Range("J86").Value = 650
$ws.Range("K86").Value = 5028
$ws.Range("L86").Value = 650
$ws.Range("M86").Value = -3905
$ws.Range("N86").Value = -2896
$ws.Range("H89").Value = 1525.6
$ws.Range("I89").Value = 5028
$ws.Range("J89").Value = 650
$ws.Range("K89").Value = 25140
$ws.Range("L89").Value = 3250
$ws.Range("M89").Value = -19524
$ws.Range("N89").Value = -14482
$ws.Range("H99").Value = 2635.3
$ws.Range("I99").Value = 2598.4707
$ws.Range("K99").Value = 2598.4707
$ws.Range("M99").Value = -1100.4707
$ws.Range("H107").Value = 1150.7241
$ws.Range("I107").Value = 1095.0769
$ws.Range("J107").Value = 1633
$ws.Range("K107").Value = 1095.0769
$ws.Range("L107").Value = 1633
$ws.Range("M107").Value = 824.9231
$ws.Range("N107").Value = -5473

# ----- CRP (31 cell updates) -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3468.7144
$ws.Range("I31").Value = 1926.5555
$ws.Range("K31").Value = 1926.5555
$ws.Range("M31").Value = -1631.5555
$ws.Range("H34").Value = 3468.7144
$ws.Range("I34").Value = 1926.5555
$ws.Range("K34").Value = 1926.5555
$ws.Range("M34").Value = -1724.5555
$ws.Range("H58").Value = 2455.0833
$ws.Range("I58").Value = 1911.1154
$ws.Range("K58").Value = 1911.1154
$ws.Range("M58").Value = -1708.1154
$ws.Range("H99").Value = 2007.8334
$ws.Range("I99").Value = 2007.8334
$ws.Range("K99").Value = 2007.8334
$ws.Range("M99").Value = -509.8334
$ws.Range("H107").Value = 28248.162
$ws.Range("I107").Value = 37675.926
$ws.Range("J107").Value = 2793.2
$ws.Range("K107").Value = 37675.926
$ws.Range("L107").Value = 2793.2
$ws.Range("M107").Value = -35755.926
$ws.Range("N107").Value = -6633.2
$ws.Range("H126").Value = 2007.8334
$ws.Range("I126").Value = 2007.8334
$ws.Range("K126").Value = 6023.5002
$ws.Range("M126").Value = -3553.5002
$ws.Range("H136").Value = 2455.0833
$ws.Range("I136").Value = 1911.1154
$ws.Range("K136").Value = 5733.3462
$ws.Range("M136").Value = -3183.3462

# ----- CUL (18 cell updates) -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41
$ws.Range("I12").Value = 25.333334
$ws.Range("J12").Value = 50.4
$ws.Range("K12").Value = 76.00000199999999
$ws.Range("L12").Value = 151.2
$ws.Range("M12").Value = 96.99999800000001
$ws.Range("N12").Value = -497.2
$ws.Range("H33").Value = 1638.7778
$ws.Range("J33").Value = 7000
$ws.Range("L33").Value = 42000
$ws.Range("N33").Value = -42566
$ws.Range("H131").Value = 1498.6562
$ws.Range("I131").Value = 1178.3334
$ws.Range("J131").Value = 1781.2941
$ws.Range("K131").Value = 3535.0002
$ws.Range("L131").Value = 5343.8823
$ws.Range("M131").Value = 1504.9998
$ws.Range("N131").Value = -15423.8823

# ----- GSM (26 cell updates) -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3499.5
$ws.Range("J80").Value = 3499.5
$ws.Range("L80").Value = 3499.5
$ws.Range("N80").Value = -5495.5
$ws.Range("H83").Value = 3499.5
$ws.Range("J83").Value = 3499.5
$ws.Range("L83").Value = 17497.5
$ws.Range("N83").Value = -27481.5
$ws.Range("H97").Value = 741.5862
$ws.Range("I97").Value = 413
$ws.Range("J97").Value = 1774.2858
$ws.Range("K97").Value = 413
$ws.Range("L97").Value = 1774.2858
$ws.Range("M97").Value = 83
$ws.Range("N97").Value = -2766.2858
$ws.Range("H121").Value = 90793.39999999999
$ws.Range("J121").Value = 90793.39999999999
$ws.Range("L121").Value = 90793.39999999999
$ws.Range("N121").Value = -94287.39999999999
$ws.Range("H126").Value = 3646.853
$ws.Range("I126").Value = 2965.0476
$ws.Range("J126").Value = 4748.231
$ws.Range("K126").Value = 8895.1428
$ws.Range("L126").Value = 14244.693
$ws.Range("M126").Value = -6425.1428
$ws.Range("N126").Value = -19184.693

# ----- LTW (40 cell updates) -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 30005
$ws.Range("I20").Value = 30005
$ws.Range("K20").Value = 30005
$ws.Range("M20").Value = -29779
$ws.Range("H55").Value = 320.29413
$ws.Range("J55").Value = 340.35294
$ws.Range("L55").Value = 340.35294
$ws.Range("N55").Value = -686.35294
$ws.Range("H74").Value = 20550.5
$ws.Range("I74").Value = 20550.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 20550.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -19552.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 20550.5
$ws.Range("I77").Value = 20550.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 61651.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -56659.5
$ws.Range("N77").ClearContents()
$ws.Range("H92").Value = 78996.664
$ws.Range("J92").Value = 78996.664
$ws.Range("L92").Value = 78996.664
$ws.Range("N92").Value = -83988.664
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4912
$ws.Range("I136").Value = 4126
$ws.Range("J136").Value = 5540.8
$ws.Range("K136").Value = 12378
$ws.Range("L136").Value = 16622.4
$ws.Range("M136").Value = -9828
$ws.Range("N136").Value = -21722.4

# ----- WVR (16 cell updates) -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5127.222
$ws.Range("I62").Value = 3822.5
$ws.Range("K62").Value = 3822.5
$ws.Range("M62").Value = -3198.5
$ws.Range("H65").Value = 5127.222
$ws.Range("I65").Value = 3822.5
$ws.Range("K65").Value = 19112.5
$ws.Range("M65").Value = -15992.5
$ws.Range("H132").Value = 5468.357
$ws.Range("I132").Value = 4698.25
$ws.Range("K132").Value = 14094.75
$ws.Range("M132").Value = -11564.75
$ws.Range("H136").Value = 2763.25
$ws.Range("I136").Value = 1789.3334
$ws.Range("K136").Value = 5368.0002
$ws.Range("M136").Value = -2818.0002

Write-Output "Applied 228 cell updates across 8 sheets."
